$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 154 ("Ají Inferno" entries),
# pushing the existing rows 154-155 (weekly snapshot dated 2021-07-22) down
# to rows 156-157 unchanged, so the new weekly snapshot (dated 2021-09-09)
# can be written into rows 154-155.
$ws.Rows("154:155").Insert()

# New row 154 - "Ají" "Inferno" "Primera" weekly update
$ws.Cells.Item(154, 1).Value = 9
$ws.Cells.Item(154, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(154, 3).Value = "Metropolitana"
$ws.Cells.Item(154, 4).Value = 44448
$ws.Cells.Item(154, 5).Value = 13
$ws.Cells.Item(154, 6).Value = 100112021
$ws.Cells.Item(154, 7).Value = "Ají"
$ws.Cells.Item(154, 8).Value = "Inferno"
$ws.Cells.Item(154, 9).Value = "Primera"
$ws.Cells.Item(154, 10).Value = 18
$ws.Cells.Item(154, 11).Value = 43000
$ws.Cells.Item(154, 12).Value = 45000
$ws.Cells.Item(154, 13).Value = 44000
$ws.Cells.Item(154, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(154, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(154, 16).Value = 3667
$ws.Cells.Item(154, 17).Value = 12
$ws.Cells.Item(154, 18).Value = "Hortaliza"

# New row 155 - "Ají" "Inferno" "Segunda" weekly update
$ws.Cells.Item(155, 1).Value = 9
$ws.Cells.Item(155, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(155, 3).Value = "Metropolitana"
$ws.Cells.Item(155, 4).Value = 44448
$ws.Cells.Item(155, 5).Value = 13
$ws.Cells.Item(155, 6).Value = 100112021
$ws.Cells.Item(155, 7).Value = "Ají"
$ws.Cells.Item(155, 8).Value = "Inferno"
$ws.Cells.Item(155, 9).Value = "Segunda"
$ws.Cells.Item(155, 10).Value = 7
$ws.Cells.Item(155, 11).Value = 41000
$ws.Cells.Item(155, 12).Value = 41000
$ws.Cells.Item(155, 13).Value = 41000
$ws.Cells.Item(155, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(155, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(155, 16).Value = 3417
$ws.Cells.Item(155, 17).Value = 12
$ws.Cells.Item(155, 18).Value = "Hortaliza"
